# Adds 14 new example rows (E0050-E0063) to the "Example" sheet, just
# below the existing last row (row 50), mirroring the row/column layout
# and formatting of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# Each entry: ID, Title, Description, Filename, Concepts
$examples = @(
    @("E0050", "Specifying duration", "This example shows the various ways in which duration can be expressed - as a period of time (6months), as an end date, or as a fixed number of occurences, or until an event occurs. The example also show some 'bad practices' where the information is expressed directly as strings - and instead shows the benefits of providing semantic information to explicitly indicate what kind of a duration it is which is necessary for unambiguous interpretations.", "E0050.ttl", "dpv:Duration,dpv:hasDuration"),
    @("E0051", "Specifying frequency", "This example shows the various ways in which frequency can be expressed, including combining frequency with duration to express complex information such as once per day for 6 months", "E0051.ttl", "dpv:Frequency,dpv:hasFrequency"),
    @("E0052", "Specifying necessity and importance in context", "This example shows a process where email address is required to be collected, and name can be optionally collected. Note that the necessity applies to the entire process i.e. both personal data and the collect processing operation. It also provides an indication of the importance of the process - for example to indicate which processes are important for the organisation (primary importance) and which are not as important or are not crucial (secondary importance).", "E0052.ttl", "dpv:Necessity,dpv:hasNecessity"),
    @("E0053", "Specifying applicability of information", "This example show how the unavailability, or non-applicability, or unknown applicability of information can be expressed using the Applicability concepts. Note that such situations may represent risks or issues that may require additional attention e.g. where the information is unknown, further steps should be taken to determine the exact applicability.", "E0053.ttl", "dpv:Applicability,dpv:hasApplicability"),
    @("E0054", "Specifying status associated with activities ", "This example shows two processes as 'activities' with the status as ongoing and proposed. The proposed activity can be useful to get an audit or approval or indicate future plans.", "E0054.ttl", "dpv:ActivityStatus,dpv:hasActivityStatus"),
    @("E0055", "Specifying compliance status and lawfulness", "This example shows the compliance status associated with activities in terms of the organisation's policies and for the EU GDPR obligations. It shows how compliance issues and lawfulness can be documented as a status associated with a process. For GDPR, it uses the concepts from EU-GDPR extension regarding lawfulness.", "E0055.ttl", "dpv:ComplianceStatus,dpv:hasComplianceStatus,eu-gdpr:GDPRLawfulness"),
    @("E0056", "Specifying the audit status assocaited with a DPIA", "This example shows how a DPIA can be documented as an audit - including a status that indicates audit is needed, and maintaining logs for how the DPIA was approved.", "E0056.ttl", "dpv:DPIA,dpv:AuditStatus,dpv:hasAuditStatus"),
    @("E0057", "Expressing GDPR Right to Data Portability could not be fulfilled due to Identity Verification failure", "The following example describes a GDPR Article 20 Data Portability request not being fulfilled due to identity verification failure. The dpv:RequestRequiresAction concept indicates further action is required - specifically to provide identity documents.", "E0057.ttl", "dpv:RightExerciseRecord,dpv:AuditStatus,dpv:hasAuditStatus,dpv:Justification,dpv:hasJustification"),
    @("E0058", "Expressing a right exercise request is delayed due to high volume of requests", "The following example uses the justification HighVolumeOfProcesses to represent a high volume of similar processes or requests causing a delay in fulfilling the rights request. The concept dpv:hasDuration is used to indicate the duration of the delay.", "E0058.ttl", "dpv:RightNonFulfilmentNotice,dpv:RequestActionDelayed,dpv:Justification,dpv:hasJustification"),
    @("E0059", "Exercising the right to rectification with contesting accuracy of information as justification", "The following example shows the justification ContestAccuracy representing contesting the accuracy of information or process to justify why the right to rectification as per GDPR Article 16 is being exercised. The information in question is represented using dpv:hasPersonalData, with two processes indicating which data should be deleted and the correction.", "E0059.ttl", "dpv:RightExerciseActivity,dpv:Justification,dpv:hasJustification"),
    @("E0060", "Specifying the location of a process", "The following example shows the use of LOC extension to express the location of a process. It also shows how the location fixture and locality concepts are useful to indicate information such as data will be stored locally and shared to a remote cloud server.", "E0060.ttl", "dpv:Location,dpv:LocationFixture,dpv:LocationLocality,dpv:hasLocation"),
    @("E0061", "Associating justifications with right exercise non-fulfilment", "The following example represents a notice outlining a failure to complete a GDPR Data Portability request due to identity verification failure.", "E0061.ttl", "dpv:hasJustification,dpv:Justification,dpv:RightNonFulfilmentNotice,dpv:hasRight"),
    @("E0062", "Using justifications across categories", "The justification concept ComplexityOfProcess represents a reason to delay a process due to the complexity of fulfilling it. To instead use it as a justification for not fulfilling the process, we create a new justification that combines the complexity of process and non-fulfilment categories.", "E0062.ttl", "dpv:hasJustification,dpv:Justification"),
    @("E0063", "Expressing data breach notifications to data subjects are not required using a justification", "The justification RightsFreedomsImpactUnlikely represents an unlikely impact on rights and freedoms, which can be used as a justification to not provide data subjects with a notification about a data breach involving their personal data as per GDPR Article 35-3b.", "E0063.ttl", "dpv:hasJustification,dpv:Justification,risk:DataBreachReport")
)

$lastRow = 50
$startRow = $lastRow + 1

for ($i = 0; $i -lt $examples.Count; $i++) {
    $row = $startRow + $i
    $entry = $examples[$i]

    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $ws.Range("D$row").Value = $entry[3]
    $ws.Range("E$row").Value = $ws.Range("E$lastRow").Value2
    $ws.Range("F$row").Value = $ws.Range("F$lastRow").Value2
    $ws.Range("G$row").Value = $entry[4]
    $ws.Range("I$row").Value = $ws.Range("I$lastRow").Value2
    $ws.Range("J$row").Value = $ws.Range("J$lastRow").Value2
    $ws.Range("K$row").Value = $ws.Range("K$lastRow").Value2

    # Copy formatting (number format / style) from the template row so the
    # new rows look identical to the existing ones.
    $ws.Range("A" + $lastRow + ":G" + $lastRow).Copy()
    $ws.Range("A" + $row + ":G" + $row).PasteSpecial(-4122)

    $ws.Range("I$lastRow").Copy()
    $ws.Range("I$row").PasteSpecial(-4122)

    $ws.Range("J$lastRow").Copy()
    $ws.Range("J$row").PasteSpecial(-4122)

    $ws.Range("K$lastRow").Copy()
    $ws.Range("K$row").PasteSpecial(-4122)
}
